$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("S4").Value = 2.3
$ws.Range("T4").Value = 1.6
$ws.Range("W4").Value = 4.33
$ws.Range("X4").Value = 1.2
